$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 147368
$ws.Range("J57").Value = 147368
$ws.Range("L57").Value = 442104
$ws.Range("N57").Value = -443102
$ws.Range("H137").Value = 3470.9666
$ws.Range("I137").Value = 1057.9166
$ws.Range("J137").Value = 4074.2292
$ws.Range("K137").Value = 3173.7498
$ws.Range("L137").Value = 12222.6876
$ws.Range("M137").Value = -623.7498000000001
$ws.Range("N137").Value = -17322.6876
$ws.Range("H140").Value = 35051.65
$ws.Range("J140").Value = 35051.65
$ws.Range("L140").Value = 35051.65
$ws.Range("N140").Value = -45411.65

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29864.918
$ws.Range("I32").Value = 30624.982
$ws.Range("J32").Value = 24001.572
$ws.Range("K32").Value = 30624.982
$ws.Range("L32").Value = 24001.572
$ws.Range("M32").Value = -30337.982
$ws.Range("N32").Value = -24575.572
$ws.Range("H113").Value = 36767.57
$ws.Range("J113").Value = 36767.57
$ws.Range("L113").Value = 36767.57
$ws.Range("N113").Value = -45445.57
$ws.Range("H122").Value = 2189.4443
$ws.Range("I122").Value = 2293.1538
$ws.Range("J122").Value = 1919.8
$ws.Range("K122").Value = 6879.4614
$ws.Range("L122").Value = 5759.4
$ws.Range("M122").Value = -4429.4614
$ws.Range("N122").Value = -10659.4
$ws.Range("H131").Value = 49936.25
$ws.Range("J131").Value = 49936.25
$ws.Range("L131").Value = 49936.25
$ws.Range("N131").Value = -60016.25
$ws.Range("H132").Value = 20002124
$ws.Range("I132").Value = 33334766
$ws.Range("K132").Value = 100004298
$ws.Range("M132").Value = -100001768

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 34998
$ws.Range("J19").Value = 34998
$ws.Range("L19").Value = 34998
$ws.Range("N19").Value = -35344
$ws.Range("H122").Value = 40774.668
$ws.Range("J122").Value = 40774.668
$ws.Range("L122").Value = 40774.668
$ws.Range("N122").Value = -50574.668
$ws.Range("H135").Value = 21262.666
$ws.Range("J135").Value = 21262.666
$ws.Range("L135").Value = 21262.666
$ws.Range("N135").Value = -31402.666
$ws.Range("H137").Value = 33770
$ws.Range("J137").Value = 33770
$ws.Range("L137").Value = 33770
$ws.Range("N137").Value = -43970

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 68571.42999999999
$ws.Range("J68").Value = 68571.42999999999
$ws.Range("L68").Value = 68571.42999999999
$ws.Range("N68").Value = -70069.42999999999
$ws.Range("H71").Value = 68571.42999999999
$ws.Range("J71").Value = 68571.42999999999
$ws.Range("L71").Value = 205714.29
$ws.Range("N71").Value = -213202.29
$ws.Range("H74").Value = 13000
$ws.Range("J74").Value = 13000
$ws.Range("L74").Value = 13000
$ws.Range("N74").Value = -14748
$ws.Range("H77").Value = 13000
$ws.Range("J77").Value = 13000
$ws.Range("L77").Value = 39000
$ws.Range("N77").Value = -47736
$ws.Range("H94").Value = 2799
$ws.Range("J94").Value = 2799
$ws.Range("L94").Value = 2799
$ws.Range("N94").Value = -3701
$ws.Range("H121").Value = 29302.375
$ws.Range("J121").Value = 29302.375
$ws.Range("L121").Value = 29302.375
$ws.Range("N121").Value = -31922.375
$ws.Range("H133").Value = 14098
$ws.Range("J133").Value = 14098
$ws.Range("L133").Value = 14098
$ws.Range("N133").Value = -19158
$ws.Range("H134").Value = 540044.1
$ws.Range("I134").Value = 1001.85
$ws.Range("J134").Value = 2336851.8
$ws.Range("K134").Value = 3005.55
$ws.Range("L134").Value = 7010555.399999999
$ws.Range("M134").Value = -470.5500000000002
$ws.Range("N134").Value = -7015625.399999999
$ws.Range("H137").Value = 27379.092
$ws.Range("J137").Value = 27379.092
$ws.Range("L137").Value = 27379.092
$ws.Range("N137").Value = -37579.092
$ws.Range("H138").Value = 39800.332
$ws.Range("J138").Value = 39800.332
$ws.Range("L138").Value = 39800.332
$ws.Range("N138").Value = -50080.332
$ws.Range("H140").Value = 14326.692
$ws.Range("J140").Value = 14326.692
$ws.Range("L140").Value = 14326.692
$ws.Range("N140").Value = -24686.692

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1164.8462
$ws.Range("J4").Value = 1638.2778
$ws.Range("L4").Value = 4914.8334
$ws.Range("N4").Value = -5138.8334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4175.278
$ws.Range("I97").Value = 1299.2307
$ws.Range("J97").Value = 11653
$ws.Range("K97").Value = 1299.2307
$ws.Range("L97").Value = 11653
$ws.Range("M97").Value = -803.2307000000001
$ws.Range("N97").Value = -12645
$ws.Range("H113").Value = 1689.7333
$ws.Range("I113").Value = 1549.6
$ws.Range("J113").Value = 1970
$ws.Range("K113").Value = 1549.6
$ws.Range("L113").Value = 1970
$ws.Range("M113").Value = 620.4000000000001
$ws.Range("N113").Value = -6310
$ws.Range("H122").Value = 1333.9412
$ws.Range("I122").Value = 1456.4166
$ws.Range("J122").Value = 1040
$ws.Range("K122").Value = 4369.2498
$ws.Range("L122").Value = 3120
$ws.Range("M122").Value = -1919.2498
$ws.Range("N122").Value = -8020
$ws.Range("H137").Value = 46695
$ws.Range("J137").Value = 46695
$ws.Range("L137").Value = 46695
$ws.Range("N137").Value = -56895
$ws.Range("H138").Value = 54500
$ws.Range("J138").Value = 54500
$ws.Range("L138").Value = 54500
$ws.Range("N138").Value = -64780
$ws.Range("H139").Value = 31950
$ws.Range("J139").Value = 31950
$ws.Range("L139").Value = 31950
$ws.Range("N139").Value = -42230

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 3587.2856
$ws.Range("I9").Value = 475
$ws.Range("J9").Value = 7737
$ws.Range("K9").Value = 475
$ws.Range("L9").Value = 7737
$ws.Range("M9").Value = -251
$ws.Range("N9").Value = -8185
$ws.Range("H55").Value = 1126.8
$ws.Range("I55").Value = 1124
$ws.Range("J55").Value = 1133.3334
$ws.Range("K55").Value = 1124
$ws.Range("L55").Value = 1133.3334
$ws.Range("M55").Value = -951
$ws.Range("N55").Value = -1479.3334
$ws.Range("H100").Value = 1685
$ws.Range("I100").Value = 1685
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1685
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1144
$ws.Range("N100").Value = $null
$ws.Range("H108").Value = 48626
$ws.Range("J108").Value = 48626
$ws.Range("L108").Value = 48626
$ws.Range("N108").Value = -56306
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null
$ws.Range("H111").Value = 43970.5
$ws.Range("J111").Value = 43970.5
$ws.Range("L111").Value = 43970.5
$ws.Range("N111").Value = -52150.5
$ws.Range("H112").Value = 24533.334
$ws.Range("J112").Value = 43600
$ws.Range("L112").Value = 43600
$ws.Range("N112").Value = -46554
$ws.Range("H114").Value = 38394
$ws.Range("J114").Value = 38394
$ws.Range("L114").Value = 38394
$ws.Range("N114").Value = -47072
$ws.Range("H116").Value = 50676
$ws.Range("J116").Value = 50676
$ws.Range("L116").Value = 50676
$ws.Range("N116").Value = -59854
$ws.Range("H120").Value = 56459.332
$ws.Range("J120").Value = 56459.332
$ws.Range("L120").Value = 56459.332
$ws.Range("N120").Value = -66135.33199999999
$ws.Range("H122").Value = 2600
$ws.Range("I122").Value = 2800
$ws.Range("J122").Value = 2533.3333
$ws.Range("K122").Value = 8400
$ws.Range("L122").Value = 7599.999899999999
$ws.Range("M122").Value = -5950
$ws.Range("N122").Value = -12499.9999
$ws.Range("H133").Value = 34870.855
$ws.Range("J133").Value = 34870.855
$ws.Range("L133").Value = 34870.855
$ws.Range("N133").Value = -39930.855
$ws.Range("H137").Value = 36662.6
$ws.Range("J137").Value = 36662.6
$ws.Range("L137").Value = 36662.6
$ws.Range("N137").Value = -46862.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 4400
$ws.Range("I20").Value = 3000
$ws.Range("K20").Value = 3000
$ws.Range("M20").Value = -2760
$ws.Range("H122").Value = 1786983.6
$ws.Range("I122").Value = 7143333
$ws.Range("J122").Value = 1533.75
$ws.Range("K122").Value = 21429999
$ws.Range("L122").Value = 4601.25
$ws.Range("M122").Value = -21427549
$ws.Range("N122").Value = -9501.25
$ws.Range("H132").Value = 3571.5417
$ws.Range("I132").Value = 4133.2334
$ws.Range("J132").Value = 2635.389
$ws.Range("K132").Value = 12399.7002
$ws.Range("L132").Value = 7906.167
$ws.Range("M132").Value = -9869.700199999999
$ws.Range("N132").Value = -12966.167
$ws.Range("H133").Value = 65280.668
$ws.Range("J133").Value = 65280.668
$ws.Range("L133").Value = 65280.668
$ws.Range("N133").Value = -75400.66800000001
$ws.Range("H139").Value = 22138.334
$ws.Range("J139").Value = 22138.334
$ws.Range("L139").Value = 22138.334
$ws.Range("N139").Value = -32418.334
